$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row (31) describing the Chinese tournament server ("比赛服"),
# following the same A/B/C/D (Host / 服务器名称 / Server name / platformId) layout
# used by the rest of the table.
$ws.Range("A31").Value = "比赛服（Tournament）"
$ws.Range("B31").Value = "艾欧尼亚"
$ws.Range("C31").Value = "Ionia"
$ws.Range("D31").Value = "FORCES"
